$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Header block updates
# ------------------------------------------------------------------
$ws.Range("E11").Value = 545406
$ws.Range("C13").Value = 2

# ------------------------------------------------------------------
# Make room for the second worker's 8 rows (Alberto) below the first
# worker's block (rows 16-23), pushing the signature block down.
# ------------------------------------------------------------------
$ws.Rows("24:31").Insert()

# Clone formatting (borders/fill/number format) for the new interior
# rows from an existing interior row, and the bottom-border row from
# the previous closing row (now row 23).
$ws.Range("B17:J17").Copy()
$ws.Range("B24:J30").PasteSpecial(-4122)

$ws.Range("B23:J23").Copy()
$ws.Range("B31:J31").PasteSpecial(-4122)

$ws.Range("B17:J17").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Worker 1: DIANA PAOLA BARBOZA HERNANDEZ (rows 16-23)
# ------------------------------------------------------------------
$periods = @("2103","2102","2101","2012","2011","2010","2009","2008")
$r = 16
foreach ($p in $periods) {
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "45764686"
    $ws.Cells.Item($r, 4).Value = "DIANA PAOLA BARBOZA HERNANDEZ"
    $ws.Cells.Item($r, 5).Value = $p
    $ws.Cells.Item($r, 6).Value = 35112
    $ws.Cells.Item($r, 7).Value = 877803
    $r = $r + 1
}

# ------------------------------------------------------------------
# Worker 2: ALBERTO CABALLERO DIAZ GRANADO (rows 24-31)
# ------------------------------------------------------------------
$fvals = @(26919,35112,35112,35112,35112,35112,35112,35112)
$r = 24
$i = 0
foreach ($p in $periods) {
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "73094794"
    $ws.Cells.Item($r, 4).Value = "ALBERTO CABALLERO DIAZ GRANADO"
    $ws.Cells.Item($r, 5).Value = $p
    $ws.Cells.Item($r, 6).Value = $fvals[$i]
    $ws.Cells.Item($r, 7).Value = 12735000
    $r = $r + 1
    $i = $i + 1
}
